{"js": "// The document renders pseudo-XML markup as plain text, one tag/line per\n// paragraph. Three paragraphs contain an `<id>...</id>` marker that was\n// originally split across three runs (`<id>`, the bare id text, `</id>`),\n// each with their own run-level formatting. The edit collapses each of\n// those three runs into a single run containing the full\n// \"<id>p002v_N</id>\" text, carrying the formatting of the first\n// (opening-tag) run.\n//\n// Re-inserting the already-combined text into the found range via\n// `insertText(..., Replace)` naturally merges the matched runs into one\n// run that keeps the formatting of the range's leading run, which is\n// exactly this transformation.\nconst ids = [\"p002v_1\", \"p002v_2\", \"p002v_3\"];\n\nfor (const id of ids) {\n  const tag = `<id>${id}</id>`;\n  const results = context.document.body.search(tag, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    const target = results.items[0];\n    target.insertText(tag, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document renders pseudo-XML markup as plain text, one tag/line per\n# paragraph. Three paragraphs contain an \"<id>...</id>\" marker that was\n# originally split across three runs (the \"<id>\" opening tag, the bare id\n# text, and the \"</id>\" closing tag), each carrying its own run-level\n# formatting. This edit collapses each of those three runs into a single\n# run containing the full \"<id>p002v_N</id>\" text, carrying the\n# formatting of the first (opening-tag) run.\n#\n# Using Find & Replace (wdReplaceAll) with the already-combined text as\n# both the search target and the replacement naturally merges the found\n# runs into a single run that keeps the formatting of the first run in\n# the match - exactly this transformation.\n$d = $word.ActiveDocument\n\n$ids = @(\"p002v_1\", \"p002v_2\", \"p002v_3\")\n\nforeach ($id in $ids) {\n    $tag = \"<id>\" + $id + \"</id>\"\n    $rng = $d.Content\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $rng.Find.Execute($tag, $true, $false, $false, $false, $false, $true, 1, $false, $tag, 2) | Out-Null\n}\n"}
